$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The former last row (37) loses its special "last row" date-only format
# and becomes a regular data row (date + time format), since a new row
# (38) is appended as today's daily update.
$ws.Range("A37").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new daily update row (row 38).
$ws.Range("A38").Value = 45778
$ws.Range("A38").NumberFormat = "YYYY-MM-DD"
$ws.Range("B38").Value = 153
$ws.Range("C38").Value = 160
$ws.Range("D38").Value = 156
